$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the identification number (A2) and index (C2) values
$ws.Range("A2").Value = 664565464
$ws.Range("C2").Value = 997653

# Update the sheet selection to D11:D12 with D11 as the active cell
$ws.Range("D11:D12").Select()
